# Applies the "Problem Solving" doc edits:
#  1. Move the "_GoBack" bookmark from the end of the "...1 to 1000." paragraph
#     to the middle of the "...bank and that they ... have all crossed."
#     sentence (splitting the run after "once they ").
#  2. Set left/right margins to 1" (1440 twips = 72 points). Top/bottom stay.

$d = $word.ActiveDocument

# --- Step 1: remove the existing _GoBack bookmark, if present -----------
# (Do this BEFORE adding the new one below, since bookmark names are
#  unique -- adding first and deleting second would remove the new one.)
$hasGoBack = $d.Bookmarks.Exists("_GoBack")
if ($hasGoBack) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: insert the _GoBack bookmark at the new location ------------
# Find the split point: right after "...once they " and before
# "have all crossed."
$findRng = $d.Content
$found = $findRng.Find.Execute("bank and that they are all still there once they ",
                       $true, $false, $false, $false, $false, $true, 1,
                       $false, "", 0)

$splitPoint = $findRng.End
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Step 3: update the page margins to 1" on left/right -----------------
$d.PageSetup.LeftMargin = 72
$d.PageSetup.RightMargin = 72

$d.Save()
